# Apply "adding transformers without normalization" edit to the
# "Computer Vision" worksheet of the Literature Review workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Computer Vision")

# ---------------------------------------------------------------------
# 1. Fill in the previously-empty WHAT? / CONTRIBUTIONS? / CATEGORY cells
#    for the existing "Bilinear Attention Networks" row (row 28).
# ---------------------------------------------------------------------
$ws.Cells.Item(28, 6).Value = "In attempts to introduce a more efficient attention mechanism, this paper proposes bilinear attention networks which are features "
$ws.Cells.Item(28, 7).Value = "(1) low-rank bilinear pooling for combining questions vector and multi-channel input image. (2) Bilinear attention network which improves upon previously used bilinear models."
$ws.Cells.Item(28, 8).Value = "Architecture"

# Row 28 grows taller to fit the new wrapped text.
$ws.Rows.Item(28).RowHeight = 75

# ---------------------------------------------------------------------
# 2. Add the new "Transformers without Normalization" row (row 29).
# ---------------------------------------------------------------------
$ws.Cells.Item(29, 1).Value = "Transformers without Normalization"
$ws.Cells.Item(29, 2).Value = "Jiachen Zhu et al."
$ws.Cells.Item(29, 3).Value = 2025
$ws.Cells.Item(29, 4).Value = "CVPR"
$ws.Cells.Item(29, 5).Value = "https://jiachenzhu.github.io/DyT/"
$ws.Cells.Item(29, 6).Value = "Jiachen Zhu et al., challenges pre-existing normalization methods employed in state-of-the art attention-based deep learning models. Specifically, they demonstrate how layer normalization layers exhibit similar behaviors as tanh functions. This papers proves this theory and demonstrates that we can achieve similar performance  accross numerous deep learning models (e.g. DiT, ViT, and ConvNeXt) by replacing  LN layers with proposed Dynamic tanh (DyT)layers. They also showcase that DyT decreases LLaMA 7B inference (7.8%) and training (8.2%) times."
$ws.Cells.Item(29, 7).Value = "(1) Discussed the purpose of layer normalization and their visual Tanh-like behavior. (2) provided pseudocode implementations for DyT. (3)Evaluated the efficiency of DyT in comparison to root mean square normalization layers. (4) Provided a series of ablation studies (i.e. replacing tanh with hard tanh and sigmoid)"
$ws.Cells.Item(29, 8).Value = "Architecture"

$ws.Rows.Item(29).RowHeight = 221

# Hyperlink the PAPER cell (E29), then restore the "Hyperlink" shared
# style (blue underline + centered wrap) used by every other link cell
# in the column by copying formats from the sibling cell E28.
$ws.Hyperlinks.Add($ws.Cells.Item(29, 5), "https://jiachenzhu.github.io/DyT/")
$ws.Cells.Item(28, 5).Copy()
$ws.Cells.Item(29, 5).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Update the view state to match where the user ended up scrolled/
#    selected after the edit.
# ---------------------------------------------------------------------
$ws.Range("F31").Select() | Out-Null
